$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: the "zh-cn" and "de-de" sheets each record a
# source file's handoff -> handback round trip. This edit records that a
# new handback just happened: the status text changes, a new "Latest Target
# File" / "Latest Handback File" pair of hyperlinked columns (E/F) gets
# filled in (duplicating the existing handoff/target hyperlinks), and the
# "Latest Handback DateTime" column (G) is stamped with the handback time.
# ---------------------------------------------------------------------------

function Update-HandbackSheet($SheetName, $XlfUrl2, $XlfUrl3, $XlfName2, $XlfName3, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    $mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/2b53d2ab2cea05bd98d5b41928f63aed46cc0cd8/e2e/25453069-af91-4079-92f1-b6b81331e315.md"
    $mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2b53d2ab2cea05bd98d5b41928f63aed46cc0cd8/e2e/689e4d4b-1aee-4e1e-b332-a648c6772d7c.md"
    $cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2b53d2ab2cea05bd98d5b41928f63aed46cc0cd8/.localization-config"

    $mdName1 = "25453069-af91-4079-92f1-b6b81331e315.md"
    $mdName2 = "689e4d4b-1aee-4e1e-b332-a648c6772d7c.md"
    $cfgName = ".localization-config"

    # Status text: the row is now handed back and in sync with en-US.
    $ws.Range("B2").Value2 = "Handed back: in sync with en-US"
    $ws.Range("B3").Value2 = "Handed back: in sync with en-US"

    # New "Latest Target File" (E) / "Latest Handback File" (F) entries -
    # mirror the same source .md / .xlf files already referenced in A / C.
    $ws.Range("E2").Value2 = $mdName1
    $ws.Range("F2").Value2 = $XlfName2
    $ws.Range("E3").Value2 = $mdName2
    $ws.Range("F3").Value2 = $XlfName3

    # Stamp the handback datetime.
    $ws.Range("G2").Value2 = $HandbackDateTime
    $ws.Range("G3").Value2 = $HandbackDateTime

    # Rebuild every hyperlink on the sheet, in column order, so the
    # relationship ids line up the way Excel would renumber them after an
    # in-the-middle insertion (rId2..rId10).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("C2"), $XlfUrl2, "", "", $XlfName2)
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfUrl2, "", "", $XlfName2)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("C3"), $XlfUrl3, "", "", $XlfName3)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfUrl3, "", "", $XlfName3)
    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)
}

$zhCnXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c18b234331138fff81f0c6d8cec4a92f391ded68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/25453069-af91-4079-92f1-b6b81331e315.7fa9607878a8c849b9b75547e90b0f0d12fc0e7c.zh-cn.xlf"
$zhCnXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c18b234331138fff81f0c6d8cec4a92f391ded68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/689e4d4b-1aee-4e1e-b332-a648c6772d7c.fd7526e35b7ba60d89053f4a3e28b9e887638de5.zh-cn.xlf"
$zhCnXlfName2 = "25453069-af91-4079-92f1-b6b81331e315.7fa9607878a8c849b9b75547e90b0f0d12fc0e7c.zh-cn.xlf"
$zhCnXlfName3 = "689e4d4b-1aee-4e1e-b332-a648c6772d7c.fd7526e35b7ba60d89053f4a3e28b9e887638de5.zh-cn.xlf"

$deDeXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df8a572b238563d88f501214ddab6711d7eaf2cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/25453069-af91-4079-92f1-b6b81331e315.7fa9607878a8c849b9b75547e90b0f0d12fc0e7c.de-de.xlf"
$deDeXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df8a572b238563d88f501214ddab6711d7eaf2cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/689e4d4b-1aee-4e1e-b332-a648c6772d7c.fd7526e35b7ba60d89053f4a3e28b9e887638de5.de-de.xlf"
$deDeXlfName2 = "25453069-af91-4079-92f1-b6b81331e315.7fa9607878a8c849b9b75547e90b0f0d12fc0e7c.de-de.xlf"
$deDeXlfName3 = "689e4d4b-1aee-4e1e-b332-a648c6772d7c.fd7526e35b7ba60d89053f4a3e28b9e887638de5.de-de.xlf"

Update-HandbackSheet "zh-cn" $zhCnXlfUrl2 $zhCnXlfUrl3 $zhCnXlfName2 $zhCnXlfName3 "2016-01-15 03:04:10"
Update-HandbackSheet "de-de" $deDeXlfUrl2 $deDeXlfUrl3 $deDeXlfName2 $deDeXlfName3 "2016-01-15 03:04:26"
